$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row data for the newly appended trading-log entries (rows 16-21),
# mirroring the existing ATTEMPT/OPENED/FAILED pattern already in the sheet.

# Row 16
$ws.Range("A16").Value = '2025-09-19T19:57:19.335693'
$ws.Range("B16").Value = 'TRADING_ATTEMPT'
$ws.Range("C16").Value = 'NEAR'
$ws.Range("D16").Value = 'UNKNOWN'
$ws.Range("E16").Value = 3.123607372056906
$ws.Range("K16").Value = 'ATTEMPT'
$ws.Range("L16").Value = 'Attempting trade 1/3'

# Row 17
$ws.Range("A17").Value = '2025-09-19T19:57:21.426956'
$ws.Range("B17").Value = 'POSITION_OPENED'
$ws.Range("C17").Value = 'NEAR'
$ws.Range("D17").Value = 'UNKNOWN'
$ws.Range("E17").Value = 3.123607372056906
$ws.Range("F17").Value = 120
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 0
$ws.Range("K17").Value = 'SUCCESS'

# Row 18
$ws.Range("A18").Value = '2025-09-19T19:57:21.442593'
$ws.Range("B18").Value = 'TRADING_ATTEMPT'
$ws.Range("C18").Value = 'SUI'
$ws.Range("D18").Value = 'UNKNOWN'
$ws.Range("E18").Value = 3.655574672649196
$ws.Range("K18").Value = 'ATTEMPT'
$ws.Range("L18").Value = 'Attempting trade 2/3'

# Row 19
$ws.Range("A19").Value = '2025-09-19T19:57:23.136489'
$ws.Range("B19").Value = 'POSITION_FAILED'
$ws.Range("C19").Value = 'SUI'
$ws.Range("D19").Value = 'UNKNOWN'
$ws.Range("K19").Value = 'FAILED'
$ws.Range("L19").Value = 'Trade execution failed for trade 2'

# Row 20
$ws.Range("A20").Value = '2025-09-19T19:57:23.152748'
$ws.Range("B20").Value = 'TRADING_ATTEMPT'
$ws.Range("C20").Value = 'ADA'
$ws.Range("D20").Value = 'UNKNOWN'
$ws.Range("E20").Value = 0.8960502649311237
$ws.Range("K20").Value = 'ATTEMPT'
$ws.Range("L20").Value = 'Attempting trade 3/3'

# Row 21
$ws.Range("A21").Value = '2025-09-19T19:57:24.887649'
$ws.Range("B21").Value = 'POSITION_FAILED'
$ws.Range("C21").Value = 'ADA'
$ws.Range("D21").Value = 'UNKNOWN'
$ws.Range("K21").Value = 'FAILED'
$ws.Range("L21").Value = 'Trade execution failed for trade 3'
